$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI ligand-receptor pair metrics recomputed with new TPM values.
# Maps each changed cell address directly to its new numeric value.
$updates = @{
    "G2" = 3.748362666666667
    "H2" = 11.245088
    "I2" = 0.2072014058556041
    "J2" = 0.2072014058556041
    "M2" = 0.1591613333333333
    "N2" = 0.477484
    "O2" = 0.08581789712686431
    "P2" = 0.08581789712686431
    "Q2" = 0.5965943998435554
    "R2" = 5.369349598591999
    "S2" = 0.01778158893225789
    "T2" = 0.0177815889322579
    "G3" = 3.748362666666667
    "H3" = 11.245088
    "I3" = 0.2072014058556041
    "J3" = 0.2072014058556041
    "O3" = 0.7091405503421059
    "P3" = 0.7091405503421059
    "Q3" = 4.929849078108444
    "R3" = 44.368641702976
    "S3" = 0.1469349189801011
    "T3" = 0.1469349189801011
    "G4" = 3.748362666666667
    "H4" = 11.245088
    "I4" = 0.2072014058556041
    "J4" = 0.2072014058556041
    "O4" = 0.2050415525310298
    "P4" = 0.2050415525310298
    "Q4" = 1.425421107608889
    "R4" = 12.82878996848
    "S4" = 0.04248489794324508
    "T4" = 0.04248489794324509
    "I5" = 0.6140431114114622
    "J5" = 0.6140431114114623
    "M5" = 0.1591613333333333
    "N5" = 0.477484
    "O5" = 0.08581789712686431
    "P5" = 0.08581789712686431
    "Q5" = 1.768012528765777
    "R5" = 15.912112758892
    "S5" = 0.05269588856656855
    "T5" = 0.05269588856656855
    "I6" = 0.6140431114114622
    "J6" = 0.6140431114114623
    "O6" = 0.7091405503421059
    "P6" = 0.7091405503421059
    "S6" = 0.4354428699601033
    "T6" = 0.4354428699601034
    "I7" = 0.6140431114114622
    "J7" = 0.6140431114114623
    "O7" = 0.2050415525310298
    "P7" = 0.2050415525310298
    "S7" = 0.1259043528847903
    "T7" = 0.1259043528847903
    "G8" = 3.233763666666667
    "H8" = 9.701291000000001
    "I8" = 0.1787554827329337
    "J8" = 0.1787554827329337
    "M8" = 0.1591613333333333
    "N8" = 0.477484
    "O8" = 0.08581789712686431
    "P8" = 0.08581789712686431
    "Q8" = 0.5146901368715555
    "R8" = 4.632211231844
    "S8" = 0.01534041962803787
    "T8" = 0.01534041962803787
    "G9" = 3.233763666666667
    "H9" = 9.701291000000001
    "I9" = 0.1787554827329337
    "J9" = 0.1787554827329337
    "O9" = 0.7091405503421059
    "P9" = 0.7091405503421059
    "Q9" = 4.253048130242445
    "R9" = 38.27743317218201
    "S9" = 0.1267627614019014
    "T9" = 0.1267627614019014
    "G10" = 3.233763666666667
    "H10" = 9.701291000000001
    "I10" = 0.1787554827329337
    "J10" = 0.1787554827329337
    "O10" = 0.2050415525310298
    "P10" = 0.2050415525310298
    "S10" = 0.03665230170299441
    "T10" = 0.03665230170299442
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
